$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "SkyBox" entry (row 3) to "SkyBox-Cubed"
$ws.Range("A3").Value = "SkyBox-Cubed"

# Add new row 10 - "Skybox-Panoramic" sampler entry, following the same
# pattern as the other sampler rows (copy row 9 then rename column A)
$ws.Range("A9:H9").Copy() | Out-Null
$ws.Range("A10:H10").PasteSpecial(-4104) | Out-Null
$ws.Range("A10").Value = "Skybox-Panoramic"

# Update the active selection to match the new row added
$ws.Range("B10:H10").Select()
